# Auto-generated Excel COM-interop script
# Applies numeric "attendee count" (F column) bumps and two refreshed
# cover-image URLs (I column), matching the gh-pages data refresh at 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 326  # was 324
$ws.Range("F4").Value = 3001  # was 2999
$ws.Range("F7").Value = 2334  # was 2333
$ws.Range("F8").Value = 1715  # was 1710
$ws.Range("F9").Value = 62  # was 61
$ws.Range("F11").Value = 135  # was 132
$ws.Range("F13").Value = 20  # was 6
$ws.Range("F14").Value = 2676  # was 2674
$ws.Range("F16").Value = 1546  # was 1545
$ws.Range("F17").Value = 7147  # was 7142
$ws.Range("F19").Value = 7289  # was 7287
$ws.Range("F21").Value = 12  # was 11
$ws.Range("F22").Value = 5586  # was 5571
$ws.Range("F23").Value = 3132  # was 3131
$ws.Range("F24").Value = 3503  # was 3502
$ws.Range("F28").Value = 1930  # was 1924
$ws.Range("F29").Value = 83  # was 82
$ws.Range("F32").Value = 228  # was 227
$ws.Range("F33").Value = 494  # was 491
$ws.Range("F34").Value = 44  # was 43
$ws.Range("F35").Value = 2454  # was 2449
$ws.Range("F36").Value = 1246  # was 1243
$ws.Range("F37").Value = 2808  # was 2801
$ws.Range("F38").Value = 52  # was 46
$ws.Range("F39").Value = 24  # was 23
$ws.Range("F40").Value = 174  # was 172
$ws.Range("F41").Value = 402  # was 401
$ws.Range("F42").Value = 1108  # was 1102
$ws.Range("F43").Value = 210  # was 208
$ws.Range("F44").Value = 489  # was 487
$ws.Range("F45").Value = 539  # was 538
$ws.Range("I41").Value = "//i1.hdslb.com/bfs/openplatform/202405/Qr2Bd5W41715931423636.jpeg"

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 223  # was 221
$ws.Range("F9").Value = 35  # was 34
$ws.Range("F12").Value = 340  # was 331

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 84  # was 83

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 326  # was 324
$ws.Range("F5").Value = 3001  # was 2999
$ws.Range("F6").Value = 2334  # was 2333
$ws.Range("F7").Value = 1715  # was 1710
$ws.Range("F8").Value = 62  # was 61
$ws.Range("F10").Value = 135  # was 132
$ws.Range("F13").Value = 84  # was 83
$ws.Range("F14").Value = 2676  # was 2674
$ws.Range("F15").Value = 1546  # was 1545
$ws.Range("F16").Value = 223  # was 221
$ws.Range("F17").Value = 35  # was 34
$ws.Range("F19").Value = 7147  # was 7142
$ws.Range("F21").Value = 7289  # was 7287
$ws.Range("F23").Value = 5586  # was 5571
$ws.Range("F24").Value = 3132  # was 3131
$ws.Range("F25").Value = 3503  # was 3502
$ws.Range("F30").Value = 1930  # was 1924
$ws.Range("F35").Value = 228  # was 227
$ws.Range("F36").Value = 494  # was 491
$ws.Range("F37").Value = 44  # was 43
$ws.Range("F38").Value = 2454  # was 2449
$ws.Range("F39").Value = 1246  # was 1243
$ws.Range("F41").Value = 2809  # was 2801
$ws.Range("F42").Value = 52  # was 46
$ws.Range("F43").Value = 24  # was 23
$ws.Range("F44").Value = 174  # was 172
$ws.Range("F45").Value = 402  # was 401
$ws.Range("F46").Value = 1108  # was 1102
$ws.Range("F47").Value = 210  # was 208
$ws.Range("F48").Value = 489  # was 487
$ws.Range("F49").Value = 539  # was 538
$ws.Range("I45").Value = "//i1.hdslb.com/bfs/openplatform/202405/Qr2Bd5W41715931423636.jpeg"

